$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                (''model'',
                 AdaBoostClassifier(estimator=LogisticRegression(C=5,
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver=''liblinear''),
                                    n_estimators=10, random_state=42))])'
$ws.Range("B2").Value = 0.6476190476190476
$ws.Range("C2").Value = '{''selector'': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), ''scaler'': StandardScaler(), ''model__n_estimators'': 10, ''model__estimator__solver'': ''liblinear'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': None, ''model__estimator__C'': 5}'
$ws.Range("D2").Value = 0.5714285714285714
$ws.Range("E2").Value = '[1 0 0 1 0 0 1 1 0 1 0 0]'
$ws.Range("F2").Value = '[1 0 1 0 0 1 1 1 1 1 1 1]'
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.5682773109243697
$ws.Range("I2").Value = 0.02528520839499841
$ws.Range("J2").Value = 0.5326330532212885
$ws.Range("K2").Value = 0.05115871885747763
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'',
                                                     random_state=42))),
                (''model'',
                 AdaBoostClassifier(estimator=LogisticRegression(C=1,
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver=''saga''),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B3").Value = 0.6476190476190476
$ws.Range("C3").Value = '{''selector'': SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'', random_state=42)), ''scaler'': None, ''model__n_estimators'': 5, ''model__estimator__solver'': ''saga'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': None, ''model__estimator__C'': 1}'
$ws.Range("D3").Value = 0.625
$ws.Range("E3").Value = '[1 0 1 0 0 0 0 1 1 0 1 1]'
$ws.Range("F3").Value = '[1 1 1 1 1 0 1 1 0 1 1 1]'
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 0.6012067840834964
$ws.Range("I3").Value = 0.02447057900489374
$ws.Range("J3").Value = 0.5441617742987606
$ws.Range("K3").Value = 0.0506619745728888
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', None), (''selector'', None),
                (''model'',
                 AdaBoostClassifier(estimator=LogisticRegression(C=0.0001,
                                                                 class_weight=''balanced'',
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver=''liblinear''),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B4").Value = 0.6285714285714284
$ws.Range("C4").Value = '{''selector'': None, ''scaler'': None, ''model__n_estimators'': 5, ''model__estimator__solver'': ''liblinear'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': ''balanced'', ''model__estimator__C'': 0.0001}'
$ws.Range("D4").Value = 0.75
$ws.Range("E4").Value = '[1 0 1 1 1 1 0 1 0 1 0 1]'
$ws.Range("F4").Value = '[1 0 1 1 1 1 1 0 1 0 0 1]'
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.6036706349206349
$ws.Range("I4").Value = 0.02423087222189142
$ws.Range("J4").Value = 0.5386243386243385
$ws.Range("K4").Value = 0.06417016374807523
